$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh recalculated TPM-derived statistics for existing rows 2-13 (E:T) ---
# Columns A-D (sending/target cluster + ligand/receptor symbol) keep the same
# text in every row; only the workbook-internal shared-string ordering moved
# upstream (no visible/semantic effect), so we leave A:D untouched here.

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.8618196666666668
$ws.Range("H2").Value = 2.585459
$ws.Range("I2").Value = 0.006773656541421759
$ws.Range("J2").Value = 0.006773656541421758
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.616015666666666
$ws.Range("N2").Value = 16.848047
$ws.Range("O2").Value = 0.2860808099623356
$ws.Range("P2").Value = 0.2860808099623357
$ws.Range("Q2").Value = 4.839992749841445
$ws.Range("R2").Value = 43.559934748573
$ws.Range("S2").Value = 0.00193781314977661
$ws.Range("T2").Value = 0.00193781314977661

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.8618196666666668
$ws.Range("H3").Value = 2.585459
$ws.Range("I3").Value = 0.006773656541421759
$ws.Range("J3").Value = 0.006773656541421758
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.435526
$ws.Range("N3").Value = 25.306578
$ws.Range("O3").Value = 0.4297071542841152
$ws.Range("P3").Value = 0.4297071542841153
$ws.Range("Q3").Value = 7.269902205478
$ws.Range("R3").Value = 65.429119849302
$ws.Range("S3").Value = 0.002910688676512326
$ws.Range("T3").Value = 0.002910688676512326

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.8618196666666668
$ws.Range("H4").Value = 2.585459
$ws.Range("I4").Value = 0.006773656541421759
$ws.Range("J4").Value = 0.006773656541421758
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.036951
$ws.Range("N4").Value = 6.110853000000001
$ws.Range("O4").Value = 0.1037626364528048
$ws.Range("P4").Value = 0.1037626364528048
$ws.Range("Q4").Value = 1.755484431836334
$ws.Range("R4").Value = 15.799359886527
$ws.Range("S4").Value = 0.0007028524611637093
$ws.Range("T4").Value = 0.0007028524611637092

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.8618196666666668
$ws.Range("H5").Value = 2.585459
$ws.Range("I5").Value = 0.006773656541421759
$ws.Range("J5").Value = 0.006773656541421758
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.542379
$ws.Range("N5").Value = 10.627137
$ws.Range("O5").Value = 0.1804493993007442
$ws.Range("P5").Value = 0.1804493993007443
$ws.Range("Q5").Value = 3.052891888987001
$ws.Range("R5").Value = 27.476027000883
$ws.Range("S5").Value = 0.001222302253969113
$ws.Range("T5").Value = 0.001222302253969113

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 125.8872733333333
$ws.Range("H6").Value = 377.66182
$ws.Range("I6").Value = 0.9894380291809874
$ws.Range("J6").Value = 0.9894380291809874
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.616015666666666
$ws.Range("N6").Value = 16.848047
$ws.Range("O6").Value = 0.2860808099623356
$ws.Range("P6").Value = 0.2860808099623357
$ws.Range("Q6").Value = 706.9848992739488
$ws.Range("R6").Value = 6362.86409346554
$ws.Range("S6").Value = 0.283059232795634
$ws.Range("T6").Value = 0.2830592327956341

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 125.8872733333333
$ws.Range("H7").Value = 377.66182
$ws.Range("I7").Value = 0.9894380291809874
$ws.Range("J7").Value = 0.9894380291809874
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.435526
$ws.Range("N7").Value = 25.306578
$ws.Range("O7").Value = 0.4297071542841152
$ws.Range("P7").Value = 0.4297071542841153
$ws.Range("Q7").Value = 1061.92536727244
$ws.Range("R7").Value = 9557.32830545196
$ws.Range("S7").Value = 0.4251685998598455
$ws.Range("T7").Value = 0.4251685998598455

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 125.8872733333333
$ws.Range("H8").Value = 377.66182
$ws.Range("I8").Value = 0.9894380291809874
$ws.Range("J8").Value = 0.9894380291809874
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.036951
$ws.Range("N8").Value = 6.110853000000001
$ws.Range("O8").Value = 0.1037626364528048
$ws.Range("P8").Value = 0.1037626364528048
$ws.Range("Q8").Value = 256.4262073036067
$ws.Range("R8").Value = 2307.83586573246
$ws.Range("S8").Value = 0.1026666985144865
$ws.Range("T8").Value = 0.1026666985144865

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 125.8872733333333
$ws.Range("H9").Value = 377.66182
$ws.Range("I9").Value = 0.9894380291809874
$ws.Range("J9").Value = 0.9894380291809874
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.542379
$ws.Range("N9").Value = 10.627137
$ws.Range("O9").Value = 0.1804493993007442
$ws.Range("P9").Value = 0.1804493993007443
$ws.Range("Q9").Value = 445.9404334232601
$ws.Range("R9").Value = 4013.463900809341
$ws.Range("S9").Value = 0.1785434980110214
$ws.Range("T9").Value = 0.1785434980110215

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.4264636666666666
$ws.Range("H10").Value = 1.279391
$ws.Range("I10").Value = 0.00335188267003504
$ws.Range("J10").Value = 0.00335188267003504
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.616015666666666
$ws.Range("N10").Value = 16.848047
$ws.Range("O10").Value = 0.2860808099623356
$ws.Range("P10").Value = 0.2860808099623357
$ws.Range("Q10").Value = 2.395026633264111
$ws.Range("R10").Value = 21.555239699377
$ws.Range("S10").Value = 0.0009589093091423404
$ws.Range("T10").Value = 0.0009589093091423406

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.4264636666666666
$ws.Range("H11").Value = 1.279391
$ws.Range("I11").Value = 0.00335188267003504
$ws.Range("J11").Value = 0.00335188267003504
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.435526
$ws.Range("N11").Value = 25.306578
$ws.Range("O11").Value = 0.4297071542841152
$ws.Range("P11").Value = 0.4297071542841153
$ws.Range("Q11").Value = 3.597445348221999
$ws.Range("R11").Value = 32.377008133998
$ws.Range("S11").Value = 0.001440327963634999
$ws.Range("T11").Value = 0.001440327963634999

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.4264636666666666
$ws.Range("H12").Value = 1.279391
$ws.Range("I12").Value = 0.00335188267003504
$ws.Range("J12").Value = 0.00335188267003504
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.036951
$ws.Range("N12").Value = 6.110853000000001
$ws.Range("O12").Value = 0.1037626364528048
$ws.Range("P12").Value = 0.1037626364528048
$ws.Range("Q12").Value = 0.8686855922803334
$ws.Range("R12").Value = 7.818170330523
$ws.Range("S12").Value = 0.0003478001829233026
$ws.Range("T12").Value = 0.0003478001829233026

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4264636666666666
$ws.Range("H13").Value = 1.279391
$ws.Range("I13").Value = 0.00335188267003504
$ws.Range("J13").Value = 0.00335188267003504
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 3.542379
$ws.Range("N13").Value = 10.627137
$ws.Range("O13").Value = 0.1804493993007442
$ws.Range("P13").Value = 0.1804493993007443
$ws.Range("Q13").Value = 1.510695937063
$ws.Range("R13").Value = 13.596263433567
$ws.Range("S13").Value = 0.0006048452143343977
$ws.Range("T13").Value = 0.0006048452143343978

# --- Append new rows 14-17: "Resolving-Mac" as sending cluster (vs every target) ---
# Row 14
$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Cthrc1"
$ws.Range("C14").Value = "Fzd5"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.05552766666666667
$ws.Range("H14").Value = 0.166583
$ws.Range("I14").Value = 0.0004364316075558192
$ws.Range("J14").Value = 0.0004364316075558192
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.616015666666666
$ws.Range("N14").Value = 16.848047
$ws.Range("O14").Value = 0.2860808099623356
$ws.Range("P14").Value = 0.2860808099623357
$ws.Range("Q14").Value = 0.3118442459334444
$ws.Range("R14").Value = 2.806598213401
$ws.Range("S14").Value = 0.000124854707782733
$ws.Range("T14").Value = 0.000124854707782733

# Row 15
$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Cthrc1"
$ws.Range("C15").Value = "Fzd5"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.05552766666666667
$ws.Range("H15").Value = 0.166583
$ws.Range("I15").Value = 0.0004364316075558192
$ws.Range("J15").Value = 0.0004364316075558192
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 8.435526
$ws.Range("N15").Value = 25.306578
$ws.Range("O15").Value = 0.4297071542841152
$ws.Range("P15").Value = 0.4297071542841153
$ws.Range("Q15").Value = 0.468405075886
$ws.Range("R15").Value = 4.215645682974
$ws.Range("S15").Value = 0.0001875377841224528
$ws.Range("T15").Value = 0.0001875377841224528

# Row 16
$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Cthrc1"
$ws.Range("C16").Value = "Fzd5"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.05552766666666667
$ws.Range("H16").Value = 0.166583
$ws.Range("I16").Value = 0.0004364316075558192
$ws.Range("J16").Value = 0.0004364316075558192
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 2.036951
$ws.Range("N16").Value = 6.110853000000001
$ws.Range("O16").Value = 0.1037626364528048
$ws.Range("P16").Value = 0.1037626364528048
$ws.Range("Q16").Value = 0.1131071361443333
$ws.Range("R16").Value = 1.017964225299
$ws.Range("S16").Value = 0.00004528529423132766
$ws.Range("T16").Value = 0.00004528529423132766

# Row 17
$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Cthrc1"
$ws.Range("C17").Value = "Fzd5"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.05552766666666667
$ws.Range("H17").Value = 0.166583
$ws.Range("I17").Value = 0.0004364316075558192
$ws.Range("J17").Value = 0.0004364316075558192
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 3.542379
$ws.Range("N17").Value = 10.627137
$ws.Range("O17").Value = 0.1804493993007442
$ws.Range("P17").Value = 0.1804493993007443
$ws.Range("Q17").Value = 0.196700040319
$ws.Range("R17").Value = 1.770300362871
$ws.Range("S17").Value = 0.00007875382141930573
$ws.Range("T17").Value = 0.00007875382141930575

